# TC_41.xlsx edit: rename sheet, refresh embedded CEIC download comment,
# widen a number format, tweak a header label and a statistic value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Replace the hidden CEIC Data Manager metadata blob stored in the A1 cell
#    comment with the refreshed download snapshot.
$newCommentText = "Jx0AAB+LCAAAAAAAAAOlGdtuG8f1VxZ8aoGSuyQtRxJGG/AmhSgpCiRVWXkphrsjcaq9sDuzkviWAi1SpCmKonCK9Io+pShQ12gTILV7+ZfAkt2n/kLPXPZGUrFWNQx759zmnDPnNkP07rXvGZckYjQM9ir1mlUxSOCELg3O9yoxP6vWH1fetVHv2iHeEY6wTzgQG8AVsN1rRvcqc84Xu6Z5dXVVu2rWwujcbFhW3XwyHEycOfFxlQaM48AhlZTLfTtXxUYd1x8Sjl3MseLcq/Qn/VqHUKcLsCEO8DmJau2Y0YAw1gs45ZQwwRkRzEmnO/yOMsxu1B7X6shcg2eU7Zh6rqIrUCq4poNtyZT6xG5Y9e2qtV1t1qfWzm6juWtt1Zr1xvsJY0qIBpjxCYkuqSMBE479hWS3tpt1a6fRtLaQuZEIZGUOsNHIc8fkkjLidojnsVIeMfUBthwOVpdzpoXMHK8W9HAVDiK8mE8p90g5NcbDtuEHWpdMiI32w4g44L8HqXRIrkaRdut0MQDsdE4jvuziZWlZx4xEo4VwUjlWG3XDgLc8EvHjBZw1cSEUAGHzKCbIvAOZMXUpc+CbBjFx7TPssTxTAYlOwuiCLbBDDiGPTSHjKvBC7ELAcco4dbJN1xDoKAoXIBE2b4eeuw9SNfEGRCq5H4CLxbbtMLzItNuERPJU5fnCmfqYJ+RrcDSZh1ejwFtO4hlzIjojbredUG/EIZGQmrsTMx76oEUGQgqWgyzhDyTgKhh1iUN97B154ERmN0FKAYBaMQ/PKO+EXuwHLNFpBYpOwKIpuU4tTNdoBIcbCKeHQT9YkbIRV+QYh1fJEa7DpRNy4BZzkuNeR6wSdwGWHN86Rp6IsHKfetAf8meRgxajYjInhG8MCYVBohTui45jt5eHsT+D9JpBjl3KXRkyMzyCOIVYB71sC7pIVf6dWtau/At6pGjUC9y76RIkgu1ye9l1wK2AENjktT0cXAD0hPL5YSuxZQMGKQ/cSb+OQ5C5Cw8vJTj1Uh6G+oHjxS5RBaEfnMkQFbqpQ70TjdZAA8hxG+FgOV0uoC4zusvhY68CnXqX8QhmgYrthHHAo6WoHMjUpG/jYfEskBtg7948ZxH5fgwjyHI/DpxO6N5/N1d55zig/P4ahnGkyuH9WaT3RGWMWZeIGiOL/r35nTI2sagUuR8QPwyoc39vg5OF9u4DDGFJVt2bg6j8uje9B31dtT2R6/dmi2B+hEZXapsWY6FDZbDq9HBz/OYdKdMlZzj2YHbj0GLP00q9CkYtdrFKkweh48hLKqAtJmMGo7Hj+jUHhgcx/tWc0BcAEybSkwky8/RiAnJILzgf4OA8hhkjrSur8LT+iv44jXDAhDnpSLFSijcToaROqVHHVsVrFMtAUMUrBCwyV+jQlPiLMMLeEBxD93XY6XkJRpEh5nO9gtbmESdxspmxplxFzRLF30Ymm5QyQyS8LpMrQEkkbFFDeEaTwZCwcghp6XWwR2eRqqpJK9+EgwPLhsOk/grjSg6KyRnARQy677fJUkzn2ULDZcjWE4QKYFFI7cn40XZjy2o2YLIRayQtHhPsGT1IZk6MfnBJGPeBbdcYE0Zd+KLY2zXeIzNCoQlKF+k2VJo7z4f2kzovVWmBvkVIkQDmjXMKbWSdMMVkDPYpwZG3zBEqUwehA3S3P/73zW+ev3rx2e1HT9988cP//uNXr/75s5tnP4KP27/+7ebjXyozFTGa4plHpELT9va21XwEcZaCkHCuKUdjN3a4hJ2eyok4XSN9sZOLTq/fORi0ZT1JgQm7aimmuDMuwzhbTpQRciN5pGYSCYrEnib1Sa8L2FyLssUd7pIUqfP4uxiVL16//Oz1yz/fya0dls1a9Z2drWq98dZRDG7D9TW6dBQbFHqAIH5UtbaqjUaOeIUGjVUDSP3Ud21x2bbgtl1Pa7mbBvImolWUljTF5+YKnwJ11HiUhkB+nSBl4E8hRVK0SoXcQofoFz9585enBSrtXQ0pSgHl5BgjNjOThRR9OJ4ak9HxuNMzpr2JiJMMl6NTwr+GWO+e5lMhqIIgxt63DGju0MyMCtyIKkZ4ZhDszI0lZGIuDwvBtgmqNnqgyFUtD6IwXqgTyTFk0A2UaTXZyLGh1kic9Oda0clQG8iVrjd//3wTgzakmw206ftHHoYKGAXK4XXWfvqvV19++OrFi9vnP7/58gcFCXqf9FkA4hyyKb9Mwx5Knu43KxB0MpHOvLC+m+svGiguUkchDTiz64/lHUqvELDWhTT5P+r70PKkYOkvgK9A0HuY9a65Tmz7EJlFAOi5wNBtw+zumQJUDc/8+p/f/u7215/ffvL8zYd/uvnojzcff/L65e/fPPuDyrrbp89vf/pMV/nVRiB1ETdaNQQa8n3EMUQ2GqJ3G1998AsjCLkBI4cRy4r01Qef5oQJReVwkkmGkS5VpKjCGmmeWfAZOVVSHQp8KYsaADqihTVTCt3EwgV1sk3erwpRIu8k4hv9aTVmxAhhmvomWFIkzpjvy6dZVEs9esdq1Bsaq7QRJswwy7n+wAtnMGQkCPkAsUJS4Pp6hoxW7ncwGLVbg4xEKTGKXBKJMFQfKBkpRUvps2SVhFoOAlgY/JzYE29Ga2TrqFRyroyZ+vHlrOWK8rf5qaJAgTpxFKmBKNBv+ZN4AcNw8kR3N14+W+bm30M1q+Yn4mzd7xbxsM5hoREW0QIg8bI0aZQqU30m3nnUOHsoXJMtAVd46gR36Od6NWldwlwZmaLu9KIojDYWnwyTkA1hkoaKYmYeT2nkmaqp283OKgEkBS/9UDc/bWHYJR7h5d6yzYx7GF4+mBfOvixrn408Vzuz3NUjdUsmIP+gLwLl/33PV8HWiiIYrMQDYOkH+OTiOob7bkltlCmSUdwAYXf9Qr5PI8afiEqgvxTkNIWcqgn1id3Uv+M8UYBTcQNTH9pILd0sqJmkLlc/q4TegPq05LXQSvK7KAR8uVioEa5fLlJEazkk1zBg5iRAUZx9D9qGekcpI00FLNTSlF+8XTJ6PudlFXtnholLZlbVmZFG9ZFrbVd3CGlW63X4FzuNhiV+TkuFQ+Wg5KrkJmZyYNnPnfb/AKFu/4EnHQAA"
$null = $ws.Range("A1").Comment.Text($newCommentText)

# 3. Widen numeric format 166 (used by AA2:AJ2) from "0.000" to "###0.000"
$ws.Range("AA2:AJ2").NumberFormat = "###0.000"

# 4. Rename the "Function Description" header to "Function Information"
$ws.Range("K1").Value = "Function Information"

# 5. Update the Kurtosis statistic in U2 with the corrected figure
$ws.Range("U2").Value = 0.2499825759175085
